# The pcap-export routine that produces the "doip" (column G) and "uds"
# (column H) byte-sequence strings used to write them as plain colon
# separated hex pairs, e.g. "02:fd:00:05:...". Downstream lookup code
# expects each byte to carry a "0x" prefix instead, e.g.
# "0x02:0xfd:0x00:0x05:...". Walk the used range and rewrite every
# matching cell in columns G and H accordingly, leaving "N/A" (and any
# other non hex-byte-list text such as headers) untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    foreach ($col in @("G", "H")) {
        $cell = $ws.Range("$col$r")
        $value = $cell.Value2

        if ($null -eq $value) {
            continue
        }

        $text = [string]$value

        if ($text -eq "N/A") {
            continue
        }

        if ($text -notmatch "^[0-9a-fA-F]{2}(:[0-9a-fA-F]{2})*$") {
            continue
        }

        $bytes = $text.Split(":")
        $prefixed = ($bytes | ForEach-Object { "0x" + $_ }) -join ":"

        $cell.Value2 = $prefixed
    }
}
